$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B14 text (cap_Delta_17O definition) ---
$ws.Range("B14").Value2 = "Triple isotopic composition of dissolved oxygen versus atmospheric O2, D17O"

# --- Update B15 text (d17O definition) ---
$ws.Range("B15").Value2 = "Enrichment of oxygen-17 in dissolved oxygen (delta(17)O) in the water body by mass spectrometry"

# --- Give B15 the same "highlighted" formatting as B16 (reuses existing style) ---
$ws.Range("B16").Font.Name = "Calibri"
$ws.Range("B16").Font.Size = 12
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Update selection to match authored state ---
$ws.Range("B16").Select() | Out-Null
